$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Nesma, Dr. Nahla Nagiub, Dr. Rana Abo-Zaid, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad'
$ws.Range("G3").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G4").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Amira Sobhy'
$ws.Range("G5").Value = 'Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama'
$ws.Range("G6").Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G7").Value = 'Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Arwa Elnagar'
$ws.Range("G8").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G9").Value = 'Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed'
$ws.Range("G16").Value = 'Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G17").Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G18").Value = 'Dr. Nesma, Dr. Nahla Nagiub, Dr. Rana Abo-Zaid, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad'
$ws.Range("G19").Value = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G20").Value = 'Dr. Nesma, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G21").Value = 'Dr. Amera Ahmad Saad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama'
$ws.Range("G22").Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G23").Value = 'Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Arwa Elnagar'
$ws.Range("G24").Value = 'Dr. Dalia Mohammad Abd Al-Salam, Dr. Madeha Saeed, Dr. Marwa Mustafa, Dr. Dina Adel, Dr. Amira Ibrahim'
$ws.Range("G25").Value = 'Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Maryam Ahmad, Dr. Arwa Al-Sayed'
$ws.Range("G32").Value = 'Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon, Dr. Nardine, Dr. Wafaa Ebida'
$ws.Range("G33").Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G34").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range("G35").Value = 'Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range("G36").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Alshimaa Atef'
$ws.Range("G37").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G38").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G40").Value = 'Dr. Marina Youhanna, Dr. Nahed Mosaad, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Nourhan Osama, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G41").Value = 'Dr. Merna Mahrous, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya'
$ws.Range("G43").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G44").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G47").Value = 'Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein'
$ws.Range("G49").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon'
$ws.Range("G50").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range("G51").Value = 'Dr. Gehan Adel, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad, Administrator'
$ws.Range("G52").Value = 'Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Alshimaa Atef'
$ws.Range("G53").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G54").Value = 'Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab'
$ws.Range("G56").Value = 'Dr. Marina Youhanna, Dr. Nahed Mosaad, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Nourhan Osama, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G57").Value = 'Dr. Merna Mahrous, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya'
$ws.Range("G59").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G60").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G63").Value = 'Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein'
$ws.Range("G65").Value = 'Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Monica, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon'
$ws.Range("G66").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G67").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G68").Value = 'Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G69").Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Lamiaa Ossama'
$ws.Range("G70").Value = 'Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda'
$ws.Range("G72").Value = 'Dr. Marina Youhanna, Dr. Nahed Mosaad, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Nourhan Osama, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G73").Value = 'Dr. Nahed Mosaad, Dr. Dalia Mohammad Abd Al-Salam, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Merna Said, Dr. Arwa Al-Sayed'
$ws.Range("G75").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G76").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G80").Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Marina Atef, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Neveen Nashaat'
$ws.Range("G81").Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Marina Sorial, Dr. Wafaa Ebida'
$ws.Range("G82").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nahla Nagiub, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Gehan Adel, Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G83").Value = 'Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G84").Value = 'Dr. Nesma, Dr. Asmaa Reda, Dr. Nourhan Mahmoud, Dr. Heba Mahmoud Ali, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad'
$ws.Range("G85").Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Lamiaa Ossama'
$ws.Range("G86").Value = 'Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda'
$ws.Range("G88").Value = 'Dr. Marina Youhanna, Dr. Nahed Mosaad, Dr. Merna Mahrous, Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Merna Said, Dr. Nourhan Osama, Dr. Sara Atawia, Dr. Basma Hamed, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Mai Mustafa, Dr. Eman M. Abo-Sakaya'
$ws.Range("G89").Value = 'Dr. Nahed Mosaad, Dr. Dalia Mohammad Abd Al-Salam, Dr. Esraa Mostafa, Dr. Madeha Saeed, Dr. Merna Said, Dr. Arwa Al-Sayed'
$ws.Range("G91").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G92").Value = 'Dr. Mona Ibrahim Hussein, Dr. Alaa Ashraf'
$ws.Range("G96").Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Marina Atef, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Neveen Nashaat'
$ws.Range("G98").Value = 'Dr. Nesma, Dr. Nahla Nagiub, Dr. Rana Abo-Zaid, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad'
$ws.Range("G100").Value = 'Dr. Nourhan Mahmoud, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Range("G101").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G102").Value = 'Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda'
$ws.Range("G103").Value = 'Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Arwa Elnagar'
$ws.Range("G104").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya'
$ws.Range("G105").Value = 'Dr. Merna Mahrous, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya'
$ws.Range("G106").Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range("G107").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G108").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G112").Value = 'Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Nahla, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon'
$ws.Range("G114").Value = 'Dr. Nesma, Dr. Nahla Nagiub, Dr. Rana Abo-Zaid, Dr. Mohammad El-Tanany, Dr. Servinaz Sayed Mohammad'
$ws.Range("G116").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Amira Sobhy'
$ws.Range("G117").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Nada Gouda, Dr. Lamiaa Ossama, Dr. Kerelos Zareef'
$ws.Range("G118").Value = 'Dr. Fatma Elhady, Dr. Amera Ahmad Saad, Dr. Nada Gouda'
$ws.Range("G119").Value = 'Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Arwa Elnagar'
$ws.Range("G120").Value = 'Dr. Marina Youhanna, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Maryam Ahmad, Dr. Amira Ibrahim, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya'
$ws.Range("G121").Value = 'Dr. Merna Mahrous, Dr. Sarah Abdelmohsen, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Nadia Mostafa, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya'
$ws.Range("G122").Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range("G123").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G124").Value = 'Dr. Sarah Mahdy, Dr. Mona Ibrahim Hussein, Dr. Youstina Gamil'
$ws.Range("G128").Value = 'Dr. Youstina Magdy, Dr. Salma Hassan, Dr. Nahla, Dr. Neveen Nashaat, Dr. Yassmen Ahmad, Dr. Remon'
